$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.954.93'
$ws.Range('E2').Value = '  +2.58%  '
$ws.Range('D3').Value = '3.082.51'
$ws.Range('E3').Value = '  +4.68%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '''579.77'
$ws.Range('D6').Value = '''168.65'
$ws.Range('E6').Value = '  +5.90%  '
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('D8').Value = '3.077.90'
$ws.Range('E8').Value = '  +4.63%  '
$ws.Range('D9').Value = '''0.524'
$ws.Range('E9').Value = '  +1.06%  '
$ws.Range('D10').Value = '''6.60'
$ws.Range('E10').Value = '  -1.28%  '
$ws.Range('E11').Value = '  +2.85%  '
$ws.Range('E12').Value = '  +5.40%  '
$ws.Range('E13').Value = '  +1.81%  '
$ws.Range('D14').Value = '''36.44'
$ws.Range('E14').Value = '  +6.05%  '
$ws.Range('E15').Value = '  -0.46%  '
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').Value = '66.910.44'
$ws.Range('E16').Value = '  +2.52%  '
$ws.Range('B17').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C17').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D17').Value = '3.594.65'
$ws.Range('E17').Value = '  +4.62%  '
$ws.Range('D18').Value = '''7.19'
$ws.Range('E18').Value = '  +3.39%  '
$ws.Range('D19').Value = '3.083.20'
$ws.Range('E19').Value = '  +4.72%  '
$ws.Range('D20').Value = '''16.19'
$ws.Range('E20').Value = '  +8.54%  '
$ws.Range('D21').Value = '''466.06'
$ws.Range('E21').Value = '  +4.58%  '
$ws.Range('E22').Value = '  +4.28%  '
$ws.Range('E23').Value = '  +4.17%  '
$ws.Range('D24').Value = '''83.32'
$ws.Range('E24').Value = '  +1.30%  '
$ws.Range('D25').Value = '''2.37'
$ws.Range('E25').Value = '  +7.30%  '
$ws.Range('D26').Value = '''12.92'
$ws.Range('E26').Value = '  +6.64%  '
$ws.Range('D27').Value = '''10.12'
$ws.Range('E27').Value = '  +1.17%  '
$ws.Range('E29').Value = '  -0.56%  '
$ws.Range('E30').Value = '  +0.64%  '
$ws.Range('E31').Value = '  +3.45%  '
$ws.Range('E32').Value = '  +0.53%  '
$ws.Range('D33').Value = '''28.15'
$ws.Range('E33').Value = '  +3.49%  '
$ws.Range('E34').Value = '  +3.35%  '
$ws.Range('D35').Value = '''0.999'
$ws.Range('E35').Value = '  +0.09%  '
$ws.Range('D36').Value = '''1.01'
$ws.Range('E36').Value = '  +3.26%  '
$ws.Range('E37').Value = '  +2.54%  '
$ws.Range('E38').Value = '  +7.21%  '
$ws.Range('D39').Value = '''46.95'
$ws.Range('E39').Value = '  +6.31%  '
$ws.Range('D40').Value = '''0.319'
$ws.Range('E40').Value = '  +7.02%  '
$ws.Range('D41').Value = '''50.19'
$ws.Range('E41').Value = '  +1.41%  '
$ws.Range('D43').Value = '''8.67'
$ws.Range('E43').Value = '  +2.61%  '
$ws.Range('D44').Value = '''2.81'
$ws.Range('E44').Value = '  -0.92%  '
$ws.Range('E45').Value = '  +2.61%  '
$ws.Range('D46').Value = '''383.77'
$ws.Range('E46').Value = '  -0.42%  '
$ws.Range('D47').Value = '2.764.53'
$ws.Range('E47').Value = '  +2.16%  '
$ws.Range('D48').Value = '''134.72'
$ws.Range('E48').Value = '  +1.72%  '
$ws.Range('E49').Value = '  +0.03%  '
$ws.Range('D50').Value = '''24.72'
$ws.Range('E50').Value = '  +6.61%  '
$ws.Range('E51').Value = '  +2.53%  '
